# Daily attendance processing - 2025-11-08 14:19:45
# Reorders the "Recorded By" (column G) names so that an exact, case-sensitive
# "System" entry is moved to the front of the comma-separated list, preserving
# the relative order of the remaining names.
#
# Note: this runtime's -eq/-ne/-ceq/-cne operators behave case-insensitively,
# so exact-case comparisons are done with the .NET string .Equals() method.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $rest = @()
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) {
                    $rest += $p
                }
            }

            $newParts = @("System") + $rest
            $newValue = $newParts -join ", "

            if (-not $newValue.Equals($value)) {
                $cell.Value = $newValue
            }
        }
    }
}
